$d = $word.ActiveDocument

# 1. Name: "DHEERAJ CHAND" -> "Dheeraj Chand"
$d.Content.Find.Execute("DHEERAJ CHAND", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dheeraj Chand", 2) | Out-Null

# 2. Professional title line -> "Professional Title"
$d.Content.Find.Execute("Research, Data Analytics & Engineering Professional", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Professional Title", 2) | Out-Null

# 3. Phone/email line
$d.Content.Find.Execute("(202) 550-7110 | Dheeraj.Chand@gmail.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "202.550.7110 | dheeraj.chand@gmail.com", 2) | Out-Null

# 4. Years of experience: 20+ -> 21
$d.Content.Find.Execute("with 20+ years", $true, $false, $false, $false, $false,
                         $true, 1, $false, "with 21 years", 2) | Out-Null

# 5. Company name placeholder
$d.Content.Find.Execute("Siege Analytics, Austin, TX | 2005", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Your Company Name, Your City, ST | 2005", 2) | Out-Null

# 6. Remove the five intervening job entries (DATA PRODUCTS MANAGER through the end of
#    RESEARCH DIRECTOR's bullets), located between the PARTNER role's last bullet and the
#    "KEY ACHIEVEMENTS AND IMPACT" heading.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($t.Contains("DATA PRODUCTS MANAGER")) { $startPara = $i }
    if ($t.Contains("Created comprehensive data visualization solutions for complex research findings")) { $endPara = $i }
}
if ($startPara -ne $null -and $endPara -ne $null) {
    $rangeStart = $d.Paragraphs.Item($startPara).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endPara).Range.End
    $d.Range($rangeStart, $rangeEnd).Delete() | Out-Null
}

# 7. Strip product codenames from the Key Achievements bullets
$d.Content.Find.Execute("Conceived and deployed BALLISTA redistricting software", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Conceived and deployed redistricting software", 2) | Out-Null

$d.Content.Find.Execute("Developed DAMON boundary estimation system", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Developed boundary estimation system", 2) | Out-Null

$d.Content.Find.Execute("Created SimCrisis econometric simulation platform", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Created econometric simulation platform", 2) | Out-Null

$d.Content.Find.Execute("Built RACSO comprehensive survey operations platform", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Built comprehensive survey operations platform", 2) | Out-Null

# 8. Remove the "Data Architecture and Engineering" and "Research Impact and Recognition"
#    sections entirely (from the "Data Architecture and Engineering" heading to the end of
#    the document body content).
$startPara2 = $null
$endPara2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($t.Contains("Data Architecture and Engineering")) { $startPara2 = $i }
    if ($t.Contains("Pioneered integration of geospatial techniques into political and market research")) { $endPara2 = $i }
}
if ($startPara2 -ne $null -and $endPara2 -ne $null) {
    $rangeStart2 = $d.Paragraphs.Item($startPara2).Range.Start
    $rangeEnd2 = $d.Paragraphs.Item($endPara2).Range.End
    $d.Range($rangeStart2, $rangeEnd2).Delete() | Out-Null
}
